$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-16 Friday" "2026-01-17 Saturday"

Replace-Text "850÷2=" "463÷8="
Replace-Text "984÷9=" "759÷3="
Replace-Text "392÷5=" "886÷8="
Replace-Text "931÷8=" "120÷5="
Replace-Text "648÷6=" "606÷7="

Replace-Text "792÷9=" "965÷8="
Replace-Text "915÷7=" "977÷9="
Replace-Text "271÷4=" "964÷5="
Replace-Text "161÷9=" "449÷5="
Replace-Text "869÷6=" "378÷3="

Replace-Text "788÷5=" "311÷8="
Replace-Text "825÷3=" "729÷7="
Replace-Text "362÷3=" "108÷3="
Replace-Text "675÷4=" "746÷8="
Replace-Text "573÷4=" "256÷6="

Replace-Text "310÷4=" "139÷7="
Replace-Text "660÷9=" "321÷4="
Replace-Text "125÷7=" "297÷9="
Replace-Text "623÷6=" "995÷8="
Replace-Text "964÷3=" "554÷9="

Replace-Text "265÷4=" "928÷8="
Replace-Text "433÷5=" "743÷9="
Replace-Text "472÷5=" "330÷2="
Replace-Text "527÷7=" "193÷4="
Replace-Text "793÷3=" "833÷3="
